$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 496, shifting existing rows 496:565 down to 497:566.
$ws.Rows(496).Insert()

# Populate the newly inserted row 496 with a fresh data record (same static
# columns as the surrounding "Macroferia Regional de Talca" / Zanahoria rows,
# with new date / volume / price values).
$ws.Range("A496").Value = 5
$ws.Range("B496").Value = "Macroferia Regional de Talca"
$ws.Range("C496").Value = "Maule"
$ws.Range("D496").Value = 45127
$ws.Range("E496").Value = 7
$ws.Range("F496").Value = 100114013
$ws.Range("G496").Value = "Zanahoria"
$ws.Range("H496").Value = "Sin especificar"
$ws.Range("I496").Value = "Primera"
$ws.Range("J496").Value = 700
$ws.Range("K496").Value = 5000
$ws.Range("L496").Value = 5000
$ws.Range("M496").Value = 5000
$ws.Range("N496").Value = '$/saco 20 kilos'
$ws.Range("O496").Value = "Región de Ñuble"
$ws.Range("P496").Value = 250
$ws.Range("Q496").Value = 20
$ws.Range("R496").Value = "Hortaliza"
